# dategraph.pptx: swap to the "legend on top / full width plot area" template.
#
# - Legend moves from the right (xlLegendPositionRight) to the top
#   (xlLegendPositionTop).
# - The plot area's manual layout is widened so the chart spans (close to)
#   the full width of the slide, instead of being squeezed to make room for
#   a right-hand legend. Vertical placement (y/h) is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$chart = $shp.Chart

# --- Legend: right -> top --------------------------------------------------
$chart.HasLegend = $true
$legend = $chart.Legend
$legend.Position = -4160   # ppLegendPositionTop

# --- Plot area: stretch across (almost) the full slide width --------------
$pa = $chart.PlotArea

# The chart's graphic frame fills the whole 12192000 x 6858000 EMU slide, so
# the manual layout fractions (c:x / c:w, 0..1 of chart width) convert to
# points via the usual 12700 EMU-per-point factor.
$pa.Left = 36.863543307086616
$pa.Width = 886.2729133858268
